$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Marzo de 2020 a las 14:50"

# Row 4
$ws.Cells.Item(4, 2).Value = 164435
$ws.Cells.Item(4, 3).Value = 591
$ws.Cells.Item(4, 5).Value = 155753
$ws.Cells.Item(4, 7).Value = 19
$ws.Cells.Item(4, 8).Value = 3175

# Row 22
$ws.Cells.Item(22, 1).Value = "Noruega"
$ws.Cells.Item(22, 2).Value = 4592
$ws.Cells.Item(22, 3).Value = 147
$ws.Cells.Item(22, 4).Value = 13
$ws.Cells.Item(22, 5).Value = 4544
$ws.Cells.Item(22, 6).Value = 97
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(22, 8).Value = 35

# Row 23
$ws.Cells.Item(23, 1).Value = "Australia"
$ws.Cells.Item(23, 2).Value = 4557
$ws.Cells.Item(23, 3).Value = 97
$ws.Cells.Item(23, 4).Value = 337
$ws.Cells.Item(23, 5).Value = 4201
$ws.Cells.Item(23, 6).Value = 28
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 19

# Row 31
$ws.Cells.Item(31, 5).Value = 1953
$ws.Cells.Item(31, 7).Value = 7
$ws.Cells.Item(31, 8).Value = 72

# Row 32
$ws.Cells.Item(32, 2).Value = 2215
$ws.Cells.Item(32, 3).Value = 160
$ws.Cells.Item(32, 5).Value = 2176
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(32, 8).Value = 32

# Row 39
$ws.Cells.Item(39, 1).Value = "Arabia Saudita"
$ws.Cells.Item(39, 2).Value = 1563
$ws.Cells.Item(39, 3).Value = 110
$ws.Cells.Item(39, 4).Value = 115
$ws.Cells.Item(39, 5).Value = 1440
$ws.Cells.Item(39, 6).Value = 12
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 8

# Row 40
$ws.Cells.Item(40, 1).Value = "Indonesia"
$ws.Cells.Item(40, 2).Value = 1528
$ws.Cells.Item(40, 3).Value = 114
$ws.Cells.Item(40, 4).Value = 81
$ws.Cells.Item(40, 5).Value = 1311
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 14
$ws.Cells.Item(40, 8).Value = 136

# Row 50
$ws.Cells.Item(50, 1).Value = "Singapur"
$ws.Cells.Item(50, 2).Value = 926
$ws.Cells.Item(50, 3).Value = 47
$ws.Cells.Item(50, 4).Value = 228
$ws.Cells.Item(50, 5).Value = 695
$ws.Cells.Item(50, 6).Value = 22
$ws.Cells.Item(50, 8).Value = 3

# Row 51
$ws.Cells.Item(51, 1).Value = "Republica Dominicana"
$ws.Cells.Item(51, 2).Value = 901
$ws.Cells.Item(51, 4).Value = 4
$ws.Cells.Item(51, 5).Value = 855
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 8).Value = 42

# Row 72
$ws.Cells.Item(72, 5).Value = 382
$ws.Cells.Item(72, 7).Value = 2
$ws.Cells.Item(72, 8).Value = 12

# Row 83
$ws.Cells.Item(83, 1).Value = "Azerbaiyan"
$ws.Cells.Item(83, 2).Value = 298
$ws.Cells.Item(83, 3).Value = 25
$ws.Cells.Item(83, 4).Value = 26
$ws.Cells.Item(83, 5).Value = 267
$ws.Cells.Item(83, 6).Value = 11
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 5

# Row 84
$ws.Cells.Item(84, 1).Value = "Kuwait"
$ws.Cells.Item(84, 2).Value = 289
$ws.Cells.Item(84, 3).Value = 23
$ws.Cells.Item(84, 4).Value = 73
$ws.Cells.Item(84, 5).Value = 216
$ws.Cells.Item(84, 6).Value = 13
$ws.Cells.Item(84, 8).Value = 0

# Row 85
$ws.Cells.Item(85, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(85, 2).Value = 285
$ws.Cells.Item(85, 4).Value = 12
$ws.Cells.Item(85, 5).Value = 266
$ws.Cells.Item(85, 6).Value = 1
$ws.Cells.Item(85, 8).Value = 7

# Row 91
$ws.Cells.Item(91, 6).Value = 4

# Row 110
$ws.Cells.Item(110, 1).Value = "Martinica"
$ws.Cells.Item(110, 2).Value = 119
$ws.Cells.Item(110, 3).Value = 26
$ws.Cells.Item(110, 4).Value = 27
$ws.Cells.Item(110, 5).Value = 90
$ws.Cells.Item(110, 6).Value = 15
$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 8).Value = 2

# Row 111
$ws.Cells.Item(111, 1).Value = "Estado de Palestina"
$ws.Cells.Item(111, 2).Value = 117
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 18
$ws.Cells.Item(111, 5).Value = 98
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 8).Value = 1

# Row 112
$ws.Cells.Item(112, 1).Value = "Georgia"
$ws.Cells.Item(112, 2).Value = 110
$ws.Cells.Item(112, 3).Value = 7
$ws.Cells.Item(112, 4).Value = 21
$ws.Cells.Item(112, 5).Value = 89
$ws.Cells.Item(112, 6).Value = 6

# Row 113
$ws.Cells.Item(113, 1).Value = "Camboya"
$ws.Cells.Item(113, 2).Value = 109
$ws.Cells.Item(113, 3).Value = 2
$ws.Cells.Item(113, 4).Value = 23
$ws.Cells.Item(113, 5).Value = 86
$ws.Cells.Item(113, 6).Value = 1

# Row 114
$ws.Cells.Item(114, 1).Value = "Kirguistan"
$ws.Cells.Item(114, 3).Value = 13
$ws.Cells.Item(114, 4).Value = 3
$ws.Cells.Item(114, 5).Value = 104
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 0

# Row 115
$ws.Cells.Item(115, 1).Value = "Bolivia"
$ws.Cells.Item(115, 2).Value = 107
$ws.Cells.Item(115, 3).Value = 10
$ws.Cells.Item(115, 4).Value = 0
$ws.Cells.Item(115, 5).Value = 101
$ws.Cells.Item(115, 6).Value = 3
$ws.Cells.Item(115, 7).Value = 2
$ws.Cells.Item(115, 8).Value = 6

# Row 116
$ws.Cells.Item(116, 1).Value = "Guadalupe"
$ws.Cells.Item(116, 2).Value = 106
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 17
$ws.Cells.Item(116, 5).Value = 85
$ws.Cells.Item(116, 6).Value = 10
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 4

# Row 117
$ws.Cells.Item(117, 1).Value = "Montenegro"
$ws.Cells.Item(117, 2).Value = 105
$ws.Cells.Item(117, 3).Value = 14
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(117, 5).Value = 103
$ws.Cells.Item(117, 6).Value = 1
$ws.Cells.Item(117, 7).Value = 1
$ws.Cells.Item(117, 8).Value = 2

# Row 118
$ws.Cells.Item(118, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(118, 2).Value = 98
$ws.Cells.Item(118, 3).Value = 17
$ws.Cells.Item(118, 4).Value = 2
$ws.Cells.Item(118, 5).Value = 88
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 8).Value = 8

# Row 119
$ws.Cells.Item(119, 1).Value = "Mayotte"
$ws.Cells.Item(119, 2).Value = 94
$ws.Cells.Item(119, 3).Value = 12
$ws.Cells.Item(119, 4).Value = 10
$ws.Cells.Item(119, 5).Value = 83
$ws.Cells.Item(119, 6).Value = 3
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 1

# Row 120
$ws.Cells.Item(120, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(120, 2).Value = 85
$ws.Cells.Item(120, 3).Value = 2
$ws.Cells.Item(120, 4).Value = 1
$ws.Cells.Item(120, 5).Value = 81
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 8).Value = 3

# Row 158
$ws.Cells.Item(158, 1).Value = "Bahamas"
$ws.Cells.Item(158, 4).Value = 1
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 0

# Row 159
$ws.Cells.Item(159, 1).Value = "Birmania"
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 7).Value = 1
$ws.Cells.Item(159, 8).Value = 1

# Row 169
$ws.Cells.Item(169, 1).Value = "Suazilandia"

# Row 170
$ws.Cells.Item(170, 1).Value = "Granada"

# Row 173
$ws.Cells.Item(173, 1).Value = "Mozambique"

# Row 174
$ws.Cells.Item(174, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(174, 3).Value = 1

# Row 175
$ws.Cells.Item(175, 1).Value = "Libia"

# Row 176
$ws.Cells.Item(176, 1).Value = "Guinea-Bisau"

# Row 177
$ws.Cells.Item(177, 1).Value = "Surinam"
$ws.Cells.Item(177, 3).Value = 0

# Row 180
$ws.Cells.Item(180, 1).Value = "Republica del Chad"
$ws.Cells.Item(180, 3).Value = 2

# Row 181
$ws.Cells.Item(181, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(181, 3).Value = 0

# Row 184
$ws.Cells.Item(184, 1).Value = "San Martin (Parte Holandesa)"

# Row 185
$ws.Cells.Item(185, 1).Value = "Santa Sede"

# Row 186
$ws.Cells.Item(186, 1).Value = "Cabo Verde"
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 8).Value = 1

# Row 188
$ws.Cells.Item(188, 1).Value = "Benin"
$ws.Cells.Item(188, 4).Value = 1
$ws.Cells.Item(188, 8).Value = 0

# Row 190
$ws.Cells.Item(190, 1).Value = "Islas Turcas y Caicos"

# Row 191
$ws.Cells.Item(191, 1).Value = "Montserrat"

# Row 198
$ws.Cells.Item(198, 1).Value = "Belice"

# Row 199
$ws.Cells.Item(199, 1).Value = "Liberia"

# Row 200
$ws.Cells.Item(200, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(200, 3).Value = 1

# Row 201
$ws.Cells.Item(201, 1).Value = "Botsuana"
$ws.Cells.Item(201, 3).Value = 0

# Row 204
$ws.Cells.Item(204, 1).Value = "Timor Oriental"

# Row 205
$ws.Cells.Item(205, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(205, 3).Value = 0

# Row 206
$ws.Cells.Item(206, 1).Value = "Sierra Leona"
$ws.Cells.Item(206, 3).Value = 1
